$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.513.79"
$ws.Range("E2").Value = "  -4.59%  "
$ws.Range("D3").Value = "2.919.83"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "2.918.15"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("E10").Value = "  -5.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.59%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000218"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "3.402.27"
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.18%  "
$ws.Range("D18").Value = "2.919.83"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").Value = "57.462.30"
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "414.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.686"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0962"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.928"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("B38").Value = "Cosmos"
$ws.Range("C38").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0678"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.107"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "373.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("D44").Value = "2.686.49"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "122.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.04%  "
